# Update the people-profile extraction sheet:
#  - Row 3: swap in a different "Grant & Project Development" profile
#  - Rows 4-7: point at refreshed source/photo URLs for the same slots
#  - Row 8: swap in a new AELC profile with full contact details
#  - Rows 9-13: newly scraped profiles appended to the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: harpreet-kaur -> soren-newman -----------------------------
$ws.Cells.Item(3, 1).Value = "https://www.uidaho.edu/cals/services/grant-and-project-development/our-people/soren-newman"
$ws.Cells.Item(3, 2).ClearContents()
$ws.Cells.Item(3, 6).Value = "Grant & Project Dev"
$ws.Cells.Item(3, 9).ClearContents()

# --- Row 4: madison-mcguire -> robyn-wakefield -------------------------
$ws.Cells.Item(4, 1).Value = "https://www.uidaho.edu/cals/people/robyn-wakefield"
$ws.Cells.Item(4, 2).Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-robin-wakefield.jpg?h=1200&la=en&w=1200&rev=d2529ad48aea4cccae4246b68ea2e9e3"

# --- Row 5: savanah-nunes -> sharon-murdock ----------------------------
$ws.Cells.Item(5, 1).Value = "https://www.uidaho.edu/cals/people/sharon-murdock"
$ws.Cells.Item(5, 2).Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-sharon-murdock.jpg?h=1200&la=en&w=1200&rev=1c743421f1754058bf05b061c742eb01"
$ws.Cells.Item(5, 6).ClearContents()
$ws.Cells.Item(5, 8).ClearContents()
$ws.Cells.Item(5, 9).ClearContents()

# --- Row 6: brandi-chastain -> stephanie-bunney ------------------------
$ws.Cells.Item(6, 1).Value = "https://www.uidaho.edu/cals/people/stephanie-bunney"
$ws.Cells.Item(6, 2).Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-stephanie-bunney.jpg?h=1200&la=en&w=1200&rev=f80e978f7d8a453393d9ffcff1d9006d"

# --- Row 7: brian-kelly -> tammy-greenwalt ------------------------------
$ws.Cells.Item(7, 1).Value = "https://www.uidaho.edu/cals/people/tammy-greenwalt"
$ws.Cells.Item(7, 2).Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-tammy-greenwalt.jpg?h=1200&la=en&w=1200&rev=a23705edf2564c34ac09e376a552d5b8"

# --- Row 8: carly-schoepflin -> amanda-moore-kriwox (AELC) -------------
$ws.Cells.Item(8, 1).Value = "https://www.uidaho.edu/cals/agricultural-education-leadership-and-communications/our-people/amanda-moore-kriwox"
$ws.Cells.Item(8, 2).Value = "/-/media/uidaho-responsive/images/cals/departments/aelc/people/1200x1200-amanda-moore-kriwox.jpg?h=1200&la=en&w=1200&rev=27d7a34316b14806af02663e3b97e24f"
$ws.Cells.Item(8, 6).Value = "Program Specialist, Academic Coordinator"
$ws.Cells.Item(8, 7).Value = "Room B-64"
$ws.Cells.Item(8, 8).Value = "208-736-3624"
$ws.Cells.Item(8, 9).Value = "akriwox@uidaho.edu"

# --- Row 9 (new): sarah-swenson (AELC) ---------------------------------
$ws.Cells.Item(9, 1).Value = "https://www.uidaho.edu/cals/agricultural-education-leadership-and-communications/our-people/sarah-swenson"
$ws.Cells.Item(9, 2).Value = "/-/media/uidaho-responsive/images/cals/departments/aelc/people/1200x1200-sarah-swenson.jpg?h=1200&la=en&w=1200&rev=504be8f9e88247bba3ffd72988b8ff84"
$ws.Cells.Item(9, 6).Value = "Administrative Coordinator"
$ws.Cells.Item(9, 7).Value = "Ag Education, Room 101"
$ws.Cells.Item(9, 8).Value = "208-885-6358"
$ws.Cells.Item(9, 9).Value = "sswenson@uidaho.edu"

# --- Row 10 (new): alexander-maas (AERS) -------------------------------
$ws.Cells.Item(10, 1).Value = "https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/alexander-maas"
$ws.Cells.Item(10, 2).Value = "/-/media/uidaho-responsive/images/cals/departments/aers/people/1200x1200-alexander-maas.jpg?h=1200&la=en&w=1200&rev=19b14f042e724ed185c93ac5ba23ee4e"

# --- Row 11 (new): andres-trujillo-barrera (AERS) ----------------------
$ws.Cells.Item(11, 1).Value = "https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/andres-trujillo-barrera"
$ws.Cells.Item(11, 2).Value = "/-/media/uidaho-responsive/images/cals/departments/aers/people/1200x1200-andres-trujillo-barrera.jpg?h=1200&la=en&w=1200&rev=40e0a170ed544ec5a45f041448cd78f6"
$ws.Cells.Item(11, 6).Value = "Associate Professor & Director, Agricultural Commodity Risk Management Program"
$ws.Cells.Item(11, 7).Value = "Ag Science, Room 37"
$ws.Cells.Item(11, 8).Value = "208-885-1151"
$ws.Cells.Item(11, 9).Value = "aatrujillo@uidaho.edu"

# --- Row 12 (new): brett-wilder (AERS) ----------------------------------
$ws.Cells.Item(12, 1).Value = "https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/brett-wilder"
$ws.Cells.Item(12, 2).Value = "/-/media/uidaho-responsive/images/cals/departments/aers/people/1200x1200-brett-wilder.jpg?h=1200&la=en&w=1200&rev=99f0a4f8cd7b4eafac548c84d342b6d1"

# --- Row 13 (new): brenda-murdoch (AVFS) --------------------------------
$ws.Cells.Item(13, 1).Value = "https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/brenda-murdoch"
$ws.Cells.Item(13, 2).Value = "/-/media/uidaho-responsive/images/cals/departments/avfs/people/1200x1200-brenda-murdoch.jpg?h=1200&la=en&w=1200&rev=ac4caf11e8b8495bbda2925b84bf85df"
